$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '60.640.28'
$ws.Range('E2').Value = '  +0.47%  '

# Row 3
$ws.Range('D3').Value = '2.379.56'
$ws.Range('E3').Value = '  -1.98%  '

# Row 4
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  -0.21%  '

# Row 5
$ws.Range('D5').Value = "'541.75"
$ws.Range('E5').Value = '  +0.66%  '

# Row 6
$ws.Range('D6').Value = "'139.90"
$ws.Range('E6').Value = '  -2.06%  '

# Row 7
$ws.Range('D7').Value = "'1.00"
$ws.Range('E7').Value = '  -0.14%  '

# Row 8
$ws.Range('E8').Value = '  -4.37%  '

# Row 9
$ws.Range('D9').Value = '2.379.68'
$ws.Range('E9').Value = '  -2.32%  '

# Row 10
$ws.Range('D10').Value = "'0.105"
$ws.Range('E10').Value = '  +0.08%  '

# Row 11
$ws.Range('E11').Value = '  +0.68%  '

# Row 12
$ws.Range('D12').Value = "'5.33"
$ws.Range('E12').Value = '  +0.86%  '

# Row 13
$ws.Range('D13').Value = "'0.342"
$ws.Range('E13').Value = '  -1.37%  '

# Row 14
$ws.Range('D14').Value = "'25.35"
$ws.Range('E14').Value = '  -0.58%  '

# Row 15
$ws.Range('D15').Value = '2.808.14'
$ws.Range('E15').Value = '  -2.45%  '

# Row 16
$ws.Range('E16').Value = '  +2.08%  '

# Row 17
$ws.Range('D17').Value = '60.354.90'
$ws.Range('E17').Value = '  +0.10%  '

# Row 18
$ws.Range('D18').Value = '2.375.85'
$ws.Range('E18').Value = '  -2.93%  '

# Row 19
$ws.Range('D19').Value = "'10.58"
$ws.Range('E19').Value = '  -3.22%  '

# Row 20
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').Value = "'4.09"
$ws.Range('E20').Value = '  -0.56%  '

# Row 21
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Value = "'316.46"
$ws.Range('E21').Value = '  +1.09%  '

# Row 22
$ws.Range('D22').Value = "'6.67"
$ws.Range('E22').Value = '  -2.52%  '

# Row 23
$ws.Range('D23').Value = "'0.998"
$ws.Range('E23').Value = '  -0.21%  '

# Row 24
$ws.Range('E24').Value = '  +6.21%  '

# Row 25
$ws.Range('D25').Value = "'63.09"
$ws.Range('E25').Value = '  +0.66%  '

# Row 26
$ws.Range('D26').Value = "'0.997"
$ws.Range('E26').Value = '  -0.54%  '

# Row 27
$ws.Range('D27').Value = '2.485.07'
$ws.Range('E27').Value = '  -4.45%  '

# Row 28
$ws.Range('B28').Value = 'Aptos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D28').Value = "'7.73"
$ws.Range('E28').Value = '  +3.83%  '

# Row 29
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0₃0921'
$ws.Range('E29').Value = '  -1.81%  '

# Row 30
$ws.Range('D30').Value = "'517.61"
$ws.Range('E30').Value = '  -0.56%  '

# Row 31
$ws.Range('D31').Value = "'1.42"
$ws.Range('E31').Value = '  -2.19%  '

# Row 32
$ws.Range('D32').Value = "'7.96"
$ws.Range('E32').Value = '  -2.47%  '

# Row 33
$ws.Range('D33').Value = "'0.144"
$ws.Range('E33').Value = '  -0.79%  '

# Row 34
$ws.Range('E34').Value = '  -2.58%  '

# Row 35
$ws.Range('E35').Value = '  +0.98%  '

# Row 36
$ws.Range('D36').Value = "'0.999"
$ws.Range('E36').Value = '  -0.19%  '

# Row 37
$ws.Range('D37').Value = "'5.45"
$ws.Range('E37').Value = '  -4.26%  '

# Row 38
$ws.Range('D38').Value = "'4.63"
$ws.Range('E38').Value = '  -2.67%  '

# Row 39
$ws.Range('D39').Value = "'0.374"
$ws.Range('E39').Value = '  +1.39%  '

# Row 40
$ws.Range('D40').Value = "'18.03"
$ws.Range('E40').Value = '  -0.58%  '

# Row 41
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').Value = "'1.72"
$ws.Range('E41').Value = '  +3.73%  '

# Row 42
$ws.Range('B42').Value = 'USDe'
$ws.Range('C42').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D42').Value = "'1.00"
$ws.Range('E42').Value = '  -0.09%  '

# Row 43
$ws.Range('D43').Value = "'137.44"
$ws.Range('E43').Value = '  -4.34%  '

# Row 44
$ws.Range('D44').Value = "'40.14"
$ws.Range('E44').Value = '  +0.64%  '

# Row 45
$ws.Range('E45').Value = '  -2.46%  '

# Row 46
$ws.Range('D46').Value = "'139.69"
$ws.Range('E46').Value = '  -3.23%  '

# Row 47
$ws.Range('D47').Value = "'3.53"
$ws.Range('E47').Value = '  +0.36%  '

# Row 48
$ws.Range('D48').Value = "'20.30"
$ws.Range('E48').Value = '  -0.81%  '

# Row 49
$ws.Range('E49').Value = '  -1.49%  '

# Row 50
$ws.Range('D50').Value = "'0.575"
$ws.Range('E50').Value = '  +0.06%  '

# Row 51
$ws.Range('D51').Value = "'0.0925"
$ws.Range('E51').Value = '  -0.36%  '
